$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1: remove header styling (bold font, thin border, centered alignment) ---
$ws.Range("A1:AR1").ClearFormats()

# --- A1 header label cleared (was "Unnamed: 0") ---
$ws.Range("A1").ClearContents()

# --- Updated / corrected numeric metrics (rows 3-8) ---
$ws.Range("E3").Value = 8
$ws.Range("G3").Value = 25
$ws.Range("H3").Value = 63
$ws.Range("I3").Value = 61
$ws.Range("J3").Value = 19
$ws.Range("K3").Value = 27
$ws.Range("L3").Value = 4
$ws.Range("M3").Value = 0
$ws.Range("P3").Value = 3
$ws.Range("Q3").Value = 52
$ws.Range("R3").Value = 35
$ws.Range("S3").Value = 26
$ws.Range("T3").Value = 19
$ws.Range("U3").Value = 11
$ws.Range("V3").Value = 7
$ws.Range("Y3").Value = 9
$ws.Range("Z3").Value = 5
$ws.Range("AA3").Value = 6
$ws.Range("AB3").Value = 1
$ws.Range("AC3").Value = 4
$ws.Range("AE3").Value = 47
$ws.Range("AF3").Value = 19
$ws.Range("AH3").Value = 9
$ws.Range("AI3").Value = 5
$ws.Range("AJ3").Value = 1
$ws.Range("AK3").Value = 10
$ws.Range("AQ3").Value = 16
$ws.Range("AR3").Value = 15
$ws.Range("E4").Value = 21
$ws.Range("G4").Value = 35
$ws.Range("H4").Value = 427
$ws.Range("I4").Value = 278
$ws.Range("J4").Value = 35
$ws.Range("K4").Value = 42
$ws.Range("L4").Value = 5
$ws.Range("M4").Value = 1
$ws.Range("P4").Value = 6
$ws.Range("Q4").Value = 186
$ws.Range("R4").Value = 63
$ws.Range("S4").Value = 70
$ws.Range("T4").Value = 23
$ws.Range("U4").Value = 20
$ws.Range("V4").Value = 12
$ws.Range("Y4").Value = 18
$ws.Range("Z4").Value = 6
$ws.Range("AA4").Value = 8
$ws.Range("AB4").Value = 2
$ws.Range("AC4").Value = 5
$ws.Range("AE4").Value = 151
$ws.Range("AF4").Value = 24
$ws.Range("AH4").Value = 12
$ws.Range("AI4").Value = 6
$ws.Range("AJ4").Value = 2
$ws.Range("AK4").Value = 12
$ws.Range("AQ4").Value = 20
$ws.Range("AR4").Value = 22
$ws.Range("E5").Value = 8592.200000000001
$ws.Range("G5").Value = 14708.15
$ws.Range("H5").Value = 141686.78
$ws.Range("I5").Value = 94165.53
$ws.Range("J5").Value = 12354.86
$ws.Range("K5").Value = 14314.16
$ws.Range("L5").Value = 2652.6
$ws.Range("M5").Value = 517.09
$ws.Range("P5").Value = 2068.68
$ws.Range("Q5").Value = 60931.1
$ws.Range("R5").Value = 24185.04
$ws.Range("S5").Value = 25860.2
$ws.Range("T5").Value = 7540.93
$ws.Range("U5").Value = 9835.32
$ws.Range("V5").Value = 4146.25
$ws.Range("Y5").Value = 6406.44
$ws.Range("Z5").Value = 3186.71
$ws.Range("AA5").Value = 4471.23
$ws.Range("AB5").Value = 333.69
$ws.Range("AC5").Value = 1376.61
$ws.Range("AE5").Value = 59772.62
$ws.Range("AF5").Value = 10170.57
$ws.Range("AH5").Value = 4621.76
$ws.Range("AI5").Value = 2285.63
$ws.Range("AJ5").Value = 1184.54
$ws.Range("AK5").Value = 5080
$ws.Range("AQ5").Value = 6306.19
$ws.Range("AR5").Value = 6881.77
$ws.Range("E6").Value = 2.7
$ws.Range("G6").Value = 4.62
$ws.Range("H6").Value = 44.53
$ws.Range("I6").Value = 29.59
$ws.Range("J6").Value = 3.88
$ws.Range("K6").Value = 4.5
$ws.Range("L6").Value = 0.83
$ws.Range("M6").Value = 0.16
$ws.Range("P6").Value = 0.65
$ws.Range("Q6").Value = 19.15
$ws.Range("R6").Value = 7.6
$ws.Range("S6").Value = 8.130000000000001
$ws.Range("T6").Value = 2.37
$ws.Range("U6").Value = 3.09
$ws.Range("V6").Value = 1.3
$ws.Range("W6").Value = 0.24
$ws.Range("X6").Value = 0.03
$ws.Range("Y6").Value = 2.01
$ws.Range("Z6").Value = 1
$ws.Range("AA6").Value = 1.41
$ws.Range("AC6").Value = 0.43
$ws.Range("AD6").Value = 0.54
$ws.Range("AE6").Value = 18.78
$ws.Range("AF6").Value = 3.2
$ws.Range("AH6").Value = 1.45
$ws.Range("AI6").Value = 0.72
$ws.Range("AJ6").Value = 0.37
$ws.Range("AK6").Value = 1.6
$ws.Range("AL6").Value = 0.7
$ws.Range("AM6").Value = 0.23
$ws.Range("AN6").Value = 0.3
$ws.Range("AO6").Value = 0.03
$ws.Range("AP6").Value = 0.38
$ws.Range("AQ6").Value = 1.98
$ws.Range("AR6").Value = 2.16
$ws.Range("E7").Value = 409.15
$ws.Range("G7").Value = 420.23
$ws.Range("H7").Value = 331.82
$ws.Range("I7").Value = 338.72
$ws.Range("J7").Value = 353
$ws.Range("K7").Value = 340.81
$ws.Range("L7").Value = 530.52
$ws.Range("M7").Value = 517.09
$ws.Range("P7").Value = 344.78
$ws.Range("Q7").Value = 327.59
$ws.Range("R7").Value = 383.89
$ws.Range("S7").Value = 369.43
$ws.Range("T7").Value = 327.87
$ws.Range("U7").Value = 491.77
$ws.Range("V7").Value = 345.52
$ws.Range("Y7").Value = 355.91
$ws.Range("Z7").Value = 531.12
$ws.Range("AA7").Value = 558.9
$ws.Range("AB7").Value = 166.85
$ws.Range("AC7").Value = 275.32
$ws.Range("AE7").Value = 395.85
$ws.Range("AF7").Value = 423.77
$ws.Range("AH7").Value = 385.15
$ws.Range("AI7").Value = 380.94
$ws.Range("AJ7").Value = 592.27
$ws.Range("AK7").Value = 423.33
$ws.Range("AQ7").Value = 315.31
$ws.Range("AR7").Value = 312.81

# --- Cells that are cleared entirely (now invalid / NaN after cleaning) ---
$ws.Range("D3").ClearContents()
$ws.Range("F3").ClearContents()
$ws.Range("D4").ClearContents()
$ws.Range("F4").ClearContents()
$ws.Range("D5").ClearContents()
$ws.Range("F5").ClearContents()
$ws.Range("D6").ClearContents()
$ws.Range("F6").ClearContents()
$ws.Range("D7").ClearContents()
$ws.Range("F7").ClearContents()
$ws.Range("D8").ClearContents()
$ws.Range("F8").ClearContents()

# --- Remove the trailing fully-blank row 10 ---
$ws.Rows.Item(10).Delete()
